# ---------------------------------------------------------------------------
# dd update en loading_bar background
#
# 1) "De achtergronden van de games ..." paragraph becomes one bold run, and
#    the path paragraph below it is split into spell-checked runs.
# 2) The old (empty) "_GoBack" bookmark near the 8bitphotos link is replaced
#    by a new "_GoBack" bookmark that wraps the new loading-bar section
#    (adding the new bookmark automatically retires the old one, since a
#    document only ever has a single "_GoBack").
# 3) The second "page break" paragraph becomes a new bold "loading bar"
#    heading + a path paragraph + a royalty-free-site heading.
# 4) The palm-trees hyperlink's visible text is swapped for the staticflickr
#    one (the link target / relationship itself is untouched).
# 5) The footer byline is re-split into many small runs with proofErr spell
#    markers around "DuckFunt" and each contributor's name.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. "De achtergronden van de games ..." + path ------------------------

$p2 = $d.Paragraphs.Item(2)
$p3 = $d.Paragraphs.Item(3)
$rng1 = $d.Range($p2.Range.Start, $p3.Range.End)
$xml1 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
  <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t>De achtergronden van de games zijn te vinden op deze locatie:</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>*:\...\GitHub\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>DuckFunt</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>docs</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Img</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>backgrounds</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
</pkg:xmlData>
'@
$rng1.InsertXML($xml1)

# --- 2 & 3. replace the second page-break paragraph with the new loading- -
#            bar section, and wrap it (plus the following path / heading /
#            hyperlink paragraphs) in the relocated "_GoBack" bookmark.

$pBreak = $d.Paragraphs.Item(28)
$xml2 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
  <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">De </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>loading</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> bar is te vinden op deze locatie:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*:\...\</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>GitHub\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>DuckFunt</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>\docs\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Img</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>load_bar</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">De achtergrond van de </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>loading</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> bar is gevonden op een royalty free site</w:t></w:r>
</w:p>
</pkg:xmlData>
'@
$pBreak.Range.InsertXML($xml2)

# --- 4. swap the displayed hyperlink text (target/relationship unchanged) -

$d.Content.Find.Execute("https://goodstock.photos/palm-trees-on-tropical-sandy-beach/", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "https://farm8.staticflickr.com/7187/6839384738_8671712e79_b.jpg", 2)

# relocate "_GoBack": adding a new one automatically removes the stale one
# that used to sit right after the 8bitphotos link.
$loadingHeading = $d.Paragraphs.Item(28)
$lastPara = $d.Paragraphs.Item(31)
$bmRange = $d.Range($loadingHeading.Range.Start, $lastPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# trailing empty paragraph added after the hyperlink paragraph
$lastPara2 = $d.Paragraphs.Item(31)
$lastPara2.Range.InsertParagraphAfter()

# --- 5. footer byline -------------------------------------------------------

$ftr = $d.Sections.Item(1).Footers.Item(1)
$xml3 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
  <w:pPr><w:pStyle w:val="Voettekst"/></w:pPr>
  <w:r><w:t xml:space="preserve">RIO4-APO2 Projectgroep: </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>DuckF</w:t></w:r>
  <w:r><w:t>unt</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>.</w:t></w:r>
  <w:r><w:br/></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t>Santino</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t>Bonora</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t xml:space="preserve">, Tom Smits, Raoul </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t>Verschoor,Floris</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t xml:space="preserve"> van Londen, Dimitri </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t>Nazari</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t xml:space="preserve">, Henk Bertens, Kevin Mertens, Anthony </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica"/>
      <w:color w:val="373E4D"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FEFEFE"/>
    </w:rPr>
    <w:t>Carincotte</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</pkg:xmlData>
'@
$ftr.Range.InsertXML($xml3)

Write-Output "edit complete"
